# Update the traffic incident list: simulate traffic re-ran and produced a
# fresh batch of incidents (some previous "avoid nodes" rows were dropped
# and new ones generated), so the sheet now has 27 data rows (A2:D28)
# instead of 30 (A2:D31).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Vehicle breakdown", 1.366393554102729, 103.7091049793846, 158067840),
    @("Roadwork", 1.344872523694638, 103.8008012700982, 139655796),
    @("Roadwork", 1.33635140407453, 103.9795047198802, 139735890),
    @("Roadwork", 1.295590034878018, 103.8924287199076, 1143824826),
    @("Roadwork", 1.399175873917032, 103.7738428336265, 5272416387),
    @("Roadwork", 1.331731353207238, 103.8177332088875, 241595350),
    @("Roadwork", 1.395202988681009, 103.8168103182145, 148992996),
    @("Heavy Traffic", 1.343493280384666, 103.6418624407713, 148133968),
    @("Roadwork", 1.40114503416915, 103.87235047762, 3636481276),
    @("Roadwork", 1.368615011940697, 103.7162652832473, 158067818),
    @("Roadwork", 1.348217635145913, 103.9549179861231, 627795097),
    @("Roadwork", 1.34895197238415, 103.6934238369014, 567939404),
    @("Roadwork", 1.360914211975437, 103.9610319873716, 6323597084),
    @("Roadwork", 1.378417484619901, 103.9238304190541, 135080726),
    @("Roadwork", 1.342822808856306, 103.7466255047426, 1218588051),
    @("Roadwork", 1.347892723763722, 103.9543987729199, 137481616),
    @("Roadwork", 1.354145196935542, 103.7272496538722, 139672081),
    @("Roadwork", 1.389113051750533, 103.7532970440217, 158044158),
    @("Roadwork", 1.304088564346271, 103.9194949747482, 259758008),
    @("Roadwork", 1.355223572003486, 103.719416471396, 1304598696),
    @("Roadwork", 1.355151924223111, 103.9633934580024, 135078416),
    @("Roadwork", 1.301567710639741, 103.9095384931241, 245248320),
    @("Roadwork", 1.292749063475387, 103.7871906402834, 566224199),
    @("Roadwork", 1.356726426744445, 103.7155388293913, 139697312),
    @("Roadwork", 1.362412908853572, 103.9590366130738, 1735669505),
    @("Roadwork", 1.365374072663697, 103.9539133285593, 2857970597),
    @("Roadwork", 1.420425919078368, 103.7714173528471, 1110672206)
)

$lastRow = $data.Count + 1

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
}

# Remove the now-unused trailing rows (previously up to row 31).
$ws.Range("A" + ($lastRow + 1) + ":D31").Clear()
